{"js": "// Remove \"Denis Bernier, \" entirely from the author list.\nconst removeResults = context.document.body.search(\"Denis Bernier, \", { matchCase: true });\nremoveResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < removeResults.items.length; i++) {\n  removeResults.items[i].insertText(\"\", \"Replace\");\n}\nawait context.sync();\n\n// Expand the trailing \" and David Sean-Fortin\" into the new author\n// sequence \", David Sean-Fortin and Pablo Vergara\".\nconst tailResults = context.document.body.search(\" and David Sean-Fortin\", { matchCase: true });\ntailResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < tailResults.items.length; i++) {\n  tailResults.items[i].insertText(\", David Sean-Fortin and Pablo Vergara\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove \"Denis Bernier, \" entirely from the author list.\n$find1 = $d.Content.Find\n$find1.Execute(\"Denis Bernier, \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# Expand the trailing \" and David Sean-Fortin\" into the new author\n# sequence \", David Sean-Fortin and Pablo Vergara\".\n$find2 = $d.Content.Find\n$find2.Execute(\" and David Sean-Fortin\", $false, $false, $false, $false, $false, $true, 1, $false, \", David Sean-Fortin and Pablo Vergara\", 2)\n"}
